$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the new number format (0.00) to the drop_check (W) column ranges.
# This creates a 5th cellXfs entry (numFmtId=2) and assigns style index 4
# to every cell in these ranges, matching the authored diff.
$ws.Range("W2:W60").NumberFormat = "0.00"
$ws.Range("W67:W125").NumberFormat = "0.00"
$ws.Range("W140:W198").NumberFormat = "0.00"
$ws.Range("W211:W269").NumberFormat = "0.00"
$ws.Range("W283:W375").NumberFormat = "0.00"

# Recompute the drop_check values as a linear ramp from 0.01 to 0.5
# across each monster-tier block (clamped at 0.5 once the block runs past
# 58 steps), replacing the old flat per-row percentages.
$dropCheckValues = @{}
$dropCheckValues[2] = 0.01
$dropCheckValues[3] = 0.018448275862068958
$dropCheckValues[4] = 0.02689655172413792
$dropCheckValues[5] = 0.03534482758620688
$dropCheckValues[6] = 0.04379310344827584
$dropCheckValues[7] = 0.0522413793103448
$dropCheckValues[8] = 0.06068965517241376
$dropCheckValues[9] = 0.06913793103448271
$dropCheckValues[10] = 0.07758620689655167
$dropCheckValues[11] = 0.08603448275862063
$dropCheckValues[12] = 0.0944827586206896
$dropCheckValues[13] = 0.10293103448275855
$dropCheckValues[14] = 0.11137931034482751
$dropCheckValues[15] = 0.11982758620689647
$dropCheckValues[16] = 0.12827586206896543
$dropCheckValues[17] = 0.1367241379310344
$dropCheckValues[18] = 0.14517241379310336
$dropCheckValues[19] = 0.1536206896551723
$dropCheckValues[20] = 0.1620689655172413
$dropCheckValues[21] = 0.17051724137931024
$dropCheckValues[22] = 0.1789655172413792
$dropCheckValues[23] = 0.18741379310344816
$dropCheckValues[24] = 0.1958620689655171
$dropCheckValues[25] = 0.2043103448275861
$dropCheckValues[26] = 0.21275862068965504
$dropCheckValues[27] = 0.22120689655172399
$dropCheckValues[28] = 0.22965517241379296
$dropCheckValues[29] = 0.2381034482758619
$dropCheckValues[30] = 0.24655172413793086
$dropCheckValues[31] = 0.25499999999999984
$dropCheckValues[32] = 0.2634482758620688
$dropCheckValues[33] = 0.27189655172413774
$dropCheckValues[34] = 0.2803448275862067
$dropCheckValues[35] = 0.2887931034482757
$dropCheckValues[36] = 0.2972413793103446
$dropCheckValues[37] = 0.3056896551724136
$dropCheckValues[38] = 0.31413793103448256
$dropCheckValues[39] = 0.3225862068965515
$dropCheckValues[40] = 0.33103448275862046
$dropCheckValues[41] = 0.33948275862068944
$dropCheckValues[42] = 0.3479310344827584
$dropCheckValues[43] = 0.35637931034482734
$dropCheckValues[44] = 0.3648275862068963
$dropCheckValues[45] = 0.3732758620689653
$dropCheckValues[46] = 0.3817241379310342
$dropCheckValues[47] = 0.3901724137931032
$dropCheckValues[48] = 0.39862068965517217
$dropCheckValues[49] = 0.4070689655172411
$dropCheckValues[50] = 0.41551724137931006
$dropCheckValues[51] = 0.42396551724137904
$dropCheckValues[52] = 0.43241379310344796
$dropCheckValues[53] = 0.44086206896551694
$dropCheckValues[54] = 0.4493103448275859
$dropCheckValues[55] = 0.45775862068965484
$dropCheckValues[56] = 0.4662068965517238
$dropCheckValues[57] = 0.4746551724137928
$dropCheckValues[58] = 0.4831034482758617
$dropCheckValues[59] = 0.4915517241379307
$dropCheckValues[60] = 0.49999999999999967
$dropCheckValues[67] = 0.01
$dropCheckValues[68] = 0.018448275862068958
$dropCheckValues[69] = 0.02689655172413792
$dropCheckValues[70] = 0.03534482758620688
$dropCheckValues[71] = 0.04379310344827584
$dropCheckValues[72] = 0.0522413793103448
$dropCheckValues[73] = 0.06068965517241376
$dropCheckValues[74] = 0.06913793103448271
$dropCheckValues[75] = 0.07758620689655167
$dropCheckValues[76] = 0.08603448275862063
$dropCheckValues[77] = 0.0944827586206896
$dropCheckValues[78] = 0.10293103448275855
$dropCheckValues[79] = 0.11137931034482751
$dropCheckValues[80] = 0.11982758620689647
$dropCheckValues[81] = 0.12827586206896543
$dropCheckValues[82] = 0.1367241379310344
$dropCheckValues[83] = 0.14517241379310336
$dropCheckValues[84] = 0.1536206896551723
$dropCheckValues[85] = 0.1620689655172413
$dropCheckValues[86] = 0.17051724137931024
$dropCheckValues[87] = 0.1789655172413792
$dropCheckValues[88] = 0.18741379310344816
$dropCheckValues[89] = 0.1958620689655171
$dropCheckValues[90] = 0.2043103448275861
$dropCheckValues[91] = 0.21275862068965504
$dropCheckValues[92] = 0.22120689655172399
$dropCheckValues[93] = 0.22965517241379296
$dropCheckValues[94] = 0.2381034482758619
$dropCheckValues[95] = 0.24655172413793086
$dropCheckValues[96] = 0.25499999999999984
$dropCheckValues[97] = 0.2634482758620688
$dropCheckValues[98] = 0.27189655172413774
$dropCheckValues[99] = 0.2803448275862067
$dropCheckValues[100] = 0.2887931034482757
$dropCheckValues[101] = 0.2972413793103446
$dropCheckValues[102] = 0.3056896551724136
$dropCheckValues[103] = 0.31413793103448256
$dropCheckValues[104] = 0.3225862068965515
$dropCheckValues[105] = 0.33103448275862046
$dropCheckValues[106] = 0.33948275862068944
$dropCheckValues[107] = 0.3479310344827584
$dropCheckValues[108] = 0.35637931034482734
$dropCheckValues[109] = 0.3648275862068963
$dropCheckValues[110] = 0.3732758620689653
$dropCheckValues[111] = 0.3817241379310342
$dropCheckValues[112] = 0.3901724137931032
$dropCheckValues[113] = 0.39862068965517217
$dropCheckValues[114] = 0.4070689655172411
$dropCheckValues[115] = 0.41551724137931006
$dropCheckValues[116] = 0.42396551724137904
$dropCheckValues[117] = 0.43241379310344796
$dropCheckValues[118] = 0.44086206896551694
$dropCheckValues[119] = 0.4493103448275859
$dropCheckValues[120] = 0.45775862068965484
$dropCheckValues[121] = 0.4662068965517238
$dropCheckValues[122] = 0.4746551724137928
$dropCheckValues[123] = 0.4831034482758617
$dropCheckValues[124] = 0.4915517241379307
$dropCheckValues[125] = 0.49999999999999967
$dropCheckValues[140] = 0.01
$dropCheckValues[141] = 0.018448275862068958
$dropCheckValues[142] = 0.02689655172413792
$dropCheckValues[143] = 0.03534482758620688
$dropCheckValues[144] = 0.04379310344827584
$dropCheckValues[145] = 0.0522413793103448
$dropCheckValues[146] = 0.06068965517241376
$dropCheckValues[147] = 0.06913793103448271
$dropCheckValues[148] = 0.07758620689655167
$dropCheckValues[149] = 0.08603448275862063
$dropCheckValues[150] = 0.0944827586206896
$dropCheckValues[151] = 0.10293103448275855
$dropCheckValues[152] = 0.11137931034482751
$dropCheckValues[153] = 0.11982758620689647
$dropCheckValues[154] = 0.12827586206896543
$dropCheckValues[155] = 0.1367241379310344
$dropCheckValues[156] = 0.14517241379310336
$dropCheckValues[157] = 0.1536206896551723
$dropCheckValues[158] = 0.1620689655172413
$dropCheckValues[159] = 0.17051724137931024
$dropCheckValues[160] = 0.1789655172413792
$dropCheckValues[161] = 0.18741379310344816
$dropCheckValues[162] = 0.1958620689655171
$dropCheckValues[163] = 0.2043103448275861
$dropCheckValues[164] = 0.21275862068965504
$dropCheckValues[165] = 0.22120689655172399
$dropCheckValues[166] = 0.22965517241379296
$dropCheckValues[167] = 0.2381034482758619
$dropCheckValues[168] = 0.24655172413793086
$dropCheckValues[169] = 0.25499999999999984
$dropCheckValues[170] = 0.2634482758620688
$dropCheckValues[171] = 0.27189655172413774
$dropCheckValues[172] = 0.2803448275862067
$dropCheckValues[173] = 0.2887931034482757
$dropCheckValues[174] = 0.2972413793103446
$dropCheckValues[175] = 0.3056896551724136
$dropCheckValues[176] = 0.31413793103448256
$dropCheckValues[177] = 0.3225862068965515
$dropCheckValues[178] = 0.33103448275862046
$dropCheckValues[179] = 0.33948275862068944
$dropCheckValues[180] = 0.3479310344827584
$dropCheckValues[181] = 0.35637931034482734
$dropCheckValues[182] = 0.3648275862068963
$dropCheckValues[183] = 0.3732758620689653
$dropCheckValues[184] = 0.3817241379310342
$dropCheckValues[185] = 0.3901724137931032
$dropCheckValues[186] = 0.39862068965517217
$dropCheckValues[187] = 0.4070689655172411
$dropCheckValues[188] = 0.41551724137931006
$dropCheckValues[189] = 0.42396551724137904
$dropCheckValues[190] = 0.43241379310344796
$dropCheckValues[191] = 0.44086206896551694
$dropCheckValues[192] = 0.4493103448275859
$dropCheckValues[193] = 0.45775862068965484
$dropCheckValues[194] = 0.4662068965517238
$dropCheckValues[195] = 0.4746551724137928
$dropCheckValues[196] = 0.4831034482758617
$dropCheckValues[197] = 0.4915517241379307
$dropCheckValues[198] = 0.49999999999999967
$dropCheckValues[211] = 0.01
$dropCheckValues[212] = 0.018448275862068958
$dropCheckValues[213] = 0.02689655172413792
$dropCheckValues[214] = 0.03534482758620688
$dropCheckValues[215] = 0.04379310344827584
$dropCheckValues[216] = 0.0522413793103448
$dropCheckValues[217] = 0.06068965517241376
$dropCheckValues[218] = 0.06913793103448271
$dropCheckValues[219] = 0.07758620689655167
$dropCheckValues[220] = 0.08603448275862063
$dropCheckValues[221] = 0.0944827586206896
$dropCheckValues[222] = 0.10293103448275855
$dropCheckValues[223] = 0.11137931034482751
$dropCheckValues[224] = 0.11982758620689647
$dropCheckValues[225] = 0.12827586206896543
$dropCheckValues[226] = 0.1367241379310344
$dropCheckValues[227] = 0.14517241379310336
$dropCheckValues[228] = 0.1536206896551723
$dropCheckValues[229] = 0.1620689655172413
$dropCheckValues[230] = 0.17051724137931024
$dropCheckValues[231] = 0.1789655172413792
$dropCheckValues[232] = 0.18741379310344816
$dropCheckValues[233] = 0.1958620689655171
$dropCheckValues[234] = 0.2043103448275861
$dropCheckValues[235] = 0.21275862068965504
$dropCheckValues[236] = 0.22120689655172399
$dropCheckValues[237] = 0.22965517241379296
$dropCheckValues[238] = 0.2381034482758619
$dropCheckValues[239] = 0.24655172413793086
$dropCheckValues[240] = 0.25499999999999984
$dropCheckValues[241] = 0.2634482758620688
$dropCheckValues[242] = 0.27189655172413774
$dropCheckValues[243] = 0.2803448275862067
$dropCheckValues[244] = 0.2887931034482757
$dropCheckValues[245] = 0.2972413793103446
$dropCheckValues[246] = 0.3056896551724136
$dropCheckValues[247] = 0.31413793103448256
$dropCheckValues[248] = 0.3225862068965515
$dropCheckValues[249] = 0.33103448275862046
$dropCheckValues[250] = 0.33948275862068944
$dropCheckValues[251] = 0.3479310344827584
$dropCheckValues[252] = 0.35637931034482734
$dropCheckValues[253] = 0.3648275862068963
$dropCheckValues[254] = 0.3732758620689653
$dropCheckValues[255] = 0.3817241379310342
$dropCheckValues[256] = 0.3901724137931032
$dropCheckValues[257] = 0.39862068965517217
$dropCheckValues[258] = 0.4070689655172411
$dropCheckValues[259] = 0.41551724137931006
$dropCheckValues[260] = 0.42396551724137904
$dropCheckValues[261] = 0.43241379310344796
$dropCheckValues[262] = 0.44086206896551694
$dropCheckValues[263] = 0.4493103448275859
$dropCheckValues[264] = 0.45775862068965484
$dropCheckValues[265] = 0.4662068965517238
$dropCheckValues[266] = 0.4746551724137928
$dropCheckValues[267] = 0.4831034482758617
$dropCheckValues[268] = 0.4915517241379307
$dropCheckValues[269] = 0.49999999999999967
$dropCheckValues[283] = 0.01
$dropCheckValues[284] = 0.018448275862068958
$dropCheckValues[285] = 0.02689655172413792
$dropCheckValues[286] = 0.03534482758620688
$dropCheckValues[287] = 0.04379310344827584
$dropCheckValues[288] = 0.0522413793103448
$dropCheckValues[289] = 0.06068965517241376
$dropCheckValues[290] = 0.06913793103448271
$dropCheckValues[291] = 0.07758620689655167
$dropCheckValues[292] = 0.08603448275862063
$dropCheckValues[293] = 0.0944827586206896
$dropCheckValues[294] = 0.10293103448275855
$dropCheckValues[295] = 0.11137931034482751
$dropCheckValues[296] = 0.11982758620689647
$dropCheckValues[297] = 0.12827586206896543
$dropCheckValues[298] = 0.1367241379310344
$dropCheckValues[299] = 0.14517241379310336
$dropCheckValues[300] = 0.1536206896551723
$dropCheckValues[301] = 0.1620689655172413
$dropCheckValues[302] = 0.17051724137931024
$dropCheckValues[303] = 0.1789655172413792
$dropCheckValues[304] = 0.18741379310344816
$dropCheckValues[305] = 0.1958620689655171
$dropCheckValues[306] = 0.2043103448275861
$dropCheckValues[307] = 0.21275862068965504
$dropCheckValues[308] = 0.22120689655172399
$dropCheckValues[309] = 0.22965517241379296
$dropCheckValues[310] = 0.2381034482758619
$dropCheckValues[311] = 0.24655172413793086
$dropCheckValues[312] = 0.25499999999999984
$dropCheckValues[313] = 0.2634482758620688
$dropCheckValues[314] = 0.27189655172413774
$dropCheckValues[315] = 0.2803448275862067
$dropCheckValues[316] = 0.2887931034482757
$dropCheckValues[317] = 0.2972413793103446
$dropCheckValues[318] = 0.3056896551724136
$dropCheckValues[319] = 0.31413793103448256
$dropCheckValues[320] = 0.3225862068965515
$dropCheckValues[321] = 0.33103448275862046
$dropCheckValues[322] = 0.33948275862068944
$dropCheckValues[323] = 0.3479310344827584
$dropCheckValues[324] = 0.35637931034482734
$dropCheckValues[325] = 0.3648275862068963
$dropCheckValues[326] = 0.3732758620689653
$dropCheckValues[327] = 0.3817241379310342
$dropCheckValues[328] = 0.3901724137931032
$dropCheckValues[329] = 0.39862068965517217
$dropCheckValues[330] = 0.4070689655172411
$dropCheckValues[331] = 0.41551724137931006
$dropCheckValues[332] = 0.42396551724137904
$dropCheckValues[333] = 0.43241379310344796
$dropCheckValues[334] = 0.44086206896551694
$dropCheckValues[335] = 0.4493103448275859
$dropCheckValues[336] = 0.45775862068965484
$dropCheckValues[337] = 0.4662068965517238
$dropCheckValues[338] = 0.4746551724137928
$dropCheckValues[339] = 0.4831034482758617
$dropCheckValues[340] = 0.4915517241379307
$dropCheckValues[341] = 0.49999999999999967
$dropCheckValues[342] = 0.49999999999999967
$dropCheckValues[343] = 0.49999999999999967
$dropCheckValues[344] = 0.49999999999999967
$dropCheckValues[345] = 0.49999999999999967
$dropCheckValues[346] = 0.49999999999999967
$dropCheckValues[347] = 0.49999999999999967
$dropCheckValues[348] = 0.49999999999999967
$dropCheckValues[349] = 0.49999999999999967
$dropCheckValues[350] = 0.49999999999999967
$dropCheckValues[351] = 0.49999999999999967
$dropCheckValues[352] = 0.49999999999999967
$dropCheckValues[353] = 0.49999999999999967
$dropCheckValues[354] = 0.49999999999999967
$dropCheckValues[355] = 0.49999999999999967
$dropCheckValues[356] = 0.49999999999999967
$dropCheckValues[357] = 0.49999999999999967
$dropCheckValues[358] = 0.49999999999999967
$dropCheckValues[359] = 0.49999999999999967
$dropCheckValues[360] = 0.49999999999999967
$dropCheckValues[361] = 0.49999999999999967
$dropCheckValues[362] = 0.49999999999999967
$dropCheckValues[363] = 0.49999999999999967
$dropCheckValues[364] = 0.49999999999999967
$dropCheckValues[365] = 0.49999999999999967
$dropCheckValues[366] = 0.49999999999999967
$dropCheckValues[367] = 0.49999999999999967
$dropCheckValues[368] = 0.49999999999999967
$dropCheckValues[369] = 0.49999999999999967
$dropCheckValues[370] = 0.49999999999999967
$dropCheckValues[371] = 0.49999999999999967
$dropCheckValues[372] = 0.49999999999999967
$dropCheckValues[373] = 0.49999999999999967
$dropCheckValues[374] = 0.49999999999999967
$dropCheckValues[375] = 0.49999999999999967

foreach ($row in $dropCheckValues.Keys) {
    $ws.Cells.Item($row, 23).Value = $dropCheckValues[$row]
}

# Update the sheet view: clear the stale frozen/scrolled-to column and
# move the active selection down onto the freshly recalculated tail of
# the drop_check column.
$ws.Range("A1").Select() | Out-Null
$ws.Range("W341:W375").Select() | Out-Null